# "shapes now read from csv"
# Adds a "Shape" column to both the Colors sheet and the pop_names sheet,
# and moves the existing "Order" column on pop_names one slot to the right.

$wb = $excel.ActiveWorkbook

$wsNames  = $wb.Worksheets.Item("pop_names")
$wsColors = $wb.Worksheets.Item("Colors")

# ---------------------------------------------------------------------------
# 1. Colors sheet: insert a new column C ("Shape") and push the old Order
#    column (previously C) to column D.
# ---------------------------------------------------------------------------
$wsColors.Columns.Item(3).Insert()

$shapeByRow = @{
    1  = 15
    2  = 15
    3  = 15
    4  = 15
    5  = 15
    6  = 15
    7  = 16
    8  = 16
    9  = 15
    10 = 18
    11 = 17
    12 = 17
    13 = 17
    14 = 18
}

foreach ($r in 1..14) {
    $wsColors.Cells.Item($r, 3).Value = $shapeByRow[$r]
}

# ---------------------------------------------------------------------------
# 2. pop_names sheet: insert a new column E ("Shape") with a VLOOKUP against
#    the Colors sheet's new Shape column, and repoint the old Order column
#    (now F) so it looks the Order value up via the new Shape column.
# ---------------------------------------------------------------------------
$wsNames.Columns.Item(5).Insert()

$wsNames.Range("E1").Value = "Shape"

foreach ($r in 2..80) {
    $wsNames.Range("E$r").Formula = "=VLOOKUP(D$r,Colors!B:C,2,FALSE)"
    $wsNames.Range("F$r").Formula = "=VLOOKUP(E$r,Colors!C:D,2,FALSE)"
}

# ---------------------------------------------------------------------------
# 3. Fix up the ranges that don't automatically follow an inserted column:
#    the autoFilter / sortState ranges and the _FilterDatabase defined name.
# ---------------------------------------------------------------------------
$wsNames.AutoFilter.Range.AutoFilter()
$wsNames.Range("A1:F29").AutoFilter()

$wb.Names.Item("_xlnm._FilterDatabase").RefersToR1C1 = "=pop_names!R1C1:R29C6"

$wsNames.Range("A1").Select()
$wsColors.Range("A1").Select()
